$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 8; $row++) {
    $cell = $ws.Range("C$row")
    $text = $cell.Value2
    if ($text -ne $null) {
        $newText = $text.Replace("You might also like ", "")
        $newText = $newText.Replace("You might also like", "")
        $cell.Value2 = $newText
    }
}
